# Regenerate the localization-status report for archive:
#   - The "Status" value "Ready for handoff" moves on to "In Translation"
#     on every sheet that shows it.
#   - Because the new status text is shorter, Excel's column autosize
#     shrinks the affected "Status" columns to the new snug width.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# The narrower width the "Status" columns settle on once the shorter
# status text is in place (closest representable width to the new
# autofit result).
$newColumnWidth = 12.5

# --- Overview sheet: Status column per-language is E (zh-cn) / F (de-de) ---
$wsOverview = $wb.Worksheets.Item("Overview")
if ($wsOverview.Range("E2").Text -eq $oldStatus) {
    $wsOverview.Range("E2").Value = $newStatus
}
if ($wsOverview.Range("F2").Text -eq $oldStatus) {
    $wsOverview.Range("F2").Value = $newStatus
}
$wsOverview.Columns(5).ColumnWidth = $newColumnWidth
$wsOverview.Columns(6).ColumnWidth = $newColumnWidth

# --- Per-language detail sheets: Status column is C ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
if ($wsZhCn.Range("C2").Text -eq $oldStatus) {
    $wsZhCn.Range("C2").Value = $newStatus
}
$wsZhCn.Columns(3).ColumnWidth = $newColumnWidth

$wsDeDe = $wb.Worksheets.Item("de-de")
if ($wsDeDe.Range("C2").Text -eq $oldStatus) {
    $wsDeDe.Range("C2").Value = $newStatus
}
$wsDeDe.Columns(3).ColumnWidth = $newColumnWidth
